$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3359.1
$ws.Range("I53").Value = 250
$ws.Range("J53").Value = 5431.8335
$ws.Range("K53").Value = 250
$ws.Range("L53").Value = 5431.8335
$ws.Range("M53").Value = 387
$ws.Range("N53").Value = -6705.8335
$ws.Range("H62").Value = 2438
$ws.Range("I62").Value = 2438
$ws.Range("K62").Value = 2438
$ws.Range("M62").Value = -1814
$ws.Range("H65").Value = 2438
$ws.Range("I65").Value = 2438
$ws.Range("K65").Value = 12190
$ws.Range("M65").Value = -9070
$ws.Range("H70").Value = 11289.1
$ws.Range("J70").Value = 1724.75
$ws.Range("L70").Value = 5174.25
$ws.Range("N70").Value = -5714.25
$ws.Range("H73").Value = 11289.1
$ws.Range("J73").Value = 1724.75
$ws.Range("L73").Value = 5174.25
$ws.Range("N73").Value = -7046.25
$ws.Range("H112").Value = 3788907.5
$ws.Range("J112").Value = 1079
$ws.Range("L112").Value = 3237
$ws.Range("N112").Value = -5453
$ws.Range("H113").Value = 71432420
$ws.Range("I113").Value = 111113500
$ws.Range("J113").Value = 6465.4
$ws.Range("K113").Value = 111113500
$ws.Range("L113").Value = 6465.4
$ws.Range("M113").Value = -111110246
$ws.Range("N113").Value = -12973.4
$ws.Range("H125").Value = 616
$ws.Range("I125").Value = 932
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 8388
$ws.Range("L125").Value = 2700
$ws.Range("M125").Value = -5928
$ws.Range("N125").Value = -7620
$ws.Range("H138").Value = 2647.3438
$ws.Range("I138").Value = 772.125
$ws.Range("J138").Value = 3272.4167
$ws.Range("K138").Value = 2316.375
$ws.Range("L138").Value = 9817.250100000001
$ws.Range("M138").Value = 2823.625
$ws.Range("N138").Value = -20097.2501

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23057.92
$ws.Range("I32").Value = 28352.184
$ws.Range("J32").Value = 6292.75
$ws.Range("K32").Value = 28352.184
$ws.Range("L32").Value = 6292.75
$ws.Range("M32").Value = -28065.184
$ws.Range("N32").Value = -6866.75
$ws.Range("H97").Value = 1478.65
$ws.Range("I97").Value = 1337.4445
$ws.Range("J97").Value = 2749.5
$ws.Range("K97").Value = 1337.4445
$ws.Range("L97").Value = 2749.5
$ws.Range("M97").Value = -841.4445000000001
$ws.Range("N97").Value = -3741.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 714.2857
$ws.Range("I20").Value = 714.2857
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 714.2857
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -467.2857
$ws.Range("N20").Value = $null
$ws.Range("H86").Value = 1637.2632
$ws.Range("I86").Value = 1308.6364
$ws.Range("J86").Value = 2089.125
$ws.Range("K86").Value = 1308.6364
$ws.Range("L86").Value = 2089.125
$ws.Range("M86").Value = -185.6364000000001
$ws.Range("N86").Value = -4335.125
$ws.Range("H89").Value = 1637.2632
$ws.Range("I89").Value = 1308.6364
$ws.Range("J89").Value = 2089.125
$ws.Range("K89").Value = 6543.182000000001
$ws.Range("L89").Value = 10445.625
$ws.Range("M89").Value = -927.1820000000007
$ws.Range("N89").Value = -21677.625
$ws.Range("H105").Value = 2220
$ws.Range("I105").Value = 2300
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 2300
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -553
$ws.Range("N105").Value = -5694

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5237.9
$ws.Range("I94").Value = 2214.1428
$ws.Range("J94").Value = 6866.077
$ws.Range("K94").Value = 2214.1428
$ws.Range("L94").Value = 6866.077
$ws.Range("M94").Value = -1763.1428
$ws.Range("N94").Value = -7768.077
$ws.Range("H132").Value = 31717.5
$ws.Range("I132").Value = 51388.5
$ws.Range("K132").Value = 154165.5
$ws.Range("M132").Value = -151635.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 810.03845
$ws.Range("I5").Value = 1049.5
$ws.Range("J5").Value = 766.5
$ws.Range("K5").Value = 3148.5
$ws.Range("L5").Value = 2299.5
$ws.Range("M5").Value = -3036.5
$ws.Range("N5").Value = -2523.5
$ws.Range("H68").Value = 1222.8182
$ws.Range("I68").Value = 417.33334
$ws.Range("K68").Value = 1252.00002
$ws.Range("M68").Value = -441.0000199999999
$ws.Range("H71").Value = 1222.8182
$ws.Range("I71").Value = 417.33334
$ws.Range("K71").Value = 3756.00006
$ws.Range("M71").Value = 299.9999399999997
$ws.Range("H80").Value = 2975.5
$ws.Range("I80").Value = 2902
$ws.Range("K80").Value = 8706
$ws.Range("M80").Value = -7770
$ws.Range("H83").Value = 2975.5
$ws.Range("I83").Value = 2902
$ws.Range("K83").Value = 26118
$ws.Range("M83").Value = -21438
$ws.Range("H131").Value = 164752.89
$ws.Range("I131").Value = 712.375
$ws.Range("J131").Value = 189513.72
$ws.Range("K131").Value = 2137.125
$ws.Range("L131").Value = 568541.16
$ws.Range("M131").Value = 2902.875
$ws.Range("N131").Value = -578621.16
$ws.Range("H135").Value = 810.03845
$ws.Range("I135").Value = 1049.5
$ws.Range("J135").Value = 766.5
$ws.Range("K135").Value = 9445.5
$ws.Range("L135").Value = 6898.5
$ws.Range("M135").Value = -6910.5
$ws.Range("N135").Value = -11968.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4270583.5
$ws.Range("I11").Value = 4780222.5
$ws.Range("J11").Value = 2741666.8
$ws.Range("K11").Value = 4780222.5
$ws.Range("L11").Value = 2741666.8
$ws.Range("M11").Value = -4780083.5
$ws.Range("N11").Value = -2741944.8
$ws.Range("H24").Value = 151000
$ws.Range("I24").Value = 187500
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 187500
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -187327
$ws.Range("N24").Value = -5346
$ws.Range("H97").Value = 1842.75
$ws.Range("I97").Value = 1012.1177
$ws.Range("J97").Value = 3860
$ws.Range("K97").Value = 1012.1177
$ws.Range("L97").Value = 3860
$ws.Range("M97").Value = -516.1177
$ws.Range("N97").Value = -4852
$ws.Range("H122").Value = 3014.7917
$ws.Range("I122").Value = 2291.5293
$ws.Range("J122").Value = 4771.2856
$ws.Range("K122").Value = 6874.5879
$ws.Range("L122").Value = 14313.8568
$ws.Range("M122").Value = -4424.5879
$ws.Range("N122").Value = -19213.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3312.6897
$ws.Range("I7").Value = 3718.4211
$ws.Range("J7").Value = 2541.8
$ws.Range("K7").Value = 3718.4211
$ws.Range("L7").Value = 2541.8
$ws.Range("M7").Value = -3606.4211
$ws.Range("N7").Value = -2765.8
$ws.Range("H45").Value = 18000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("H93").Value = 3313.2666
$ws.Range("J93").Value = 3234
$ws.Range("L93").Value = 3234
$ws.Range("N93").Value = -5730
$ws.Range("H122").Value = 3047.2942
$ws.Range("I122").Value = 3160.8
$ws.Range("K122").Value = 9482.400000000001
$ws.Range("M122").Value = -7032.400000000001
$ws.Range("H126").Value = 3312.6897
$ws.Range("I126").Value = 3718.4211
$ws.Range("J126").Value = 2541.8
$ws.Range("K126").Value = 11155.2633
$ws.Range("L126").Value = 7625.400000000001
$ws.Range("M126").Value = -8685.263300000001
$ws.Range("N126").Value = -12565.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3227816.2
$ws.Range("I136").Value = 8065266
$ws.Range("J136").Value = 2849.8333
$ws.Range("K136").Value = 24195798
$ws.Range("L136").Value = 8549.499899999999
$ws.Range("M136").Value = -24193248
$ws.Range("N136").Value = -13649.4999

Write-Output "Applied all Typhon_Profits market data updates."